# Applies the "Updated symbol list" crypto-price refresh described by the
# commit diff: most Price (D) / Volume(1h) (E) cells get refreshed text
# values, and rows 41/42 (BKEXToken <-> KickToken) swap places along with
# their Link/Price/Volume columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain-text cells (Coin name / Link columns) ----------------------------
$textUpdates = @(
    @{ Cell = 'B41'; Value = 'KickToken' }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick' }
    @{ Cell = 'B42'; Value = 'BKEXToken' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk' }
)

foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# --- Numeric-looking text cells (Price / Volume(1h) columns) ----------------
# These must stay plain text (e.g. "243.90", "-0.12%") rather than being
# auto-converted to a number/percentage by Excel, so force a text format
# before assigning, then restore the default (un-styled) cell style so no
# stray style index gets attached to the cell.
$numericTextUpdates = @(
    @{ Cell = 'D2'; Value = '243.90' }
    @{ Cell = 'E2'; Value = '-0.12%' }
    @{ Cell = 'D3'; Value = '30.15' }
    @{ Cell = 'E3'; Value = '14.07%' }
    @{ Cell = 'E4'; Value = '-0.16%' }
    @{ Cell = 'D5'; Value = '0.05675' }
    @{ Cell = 'E5'; Value = '1.24%' }
    @{ Cell = 'D6'; Value = '6.538' }
    @{ Cell = 'E6'; Value = '1.01%' }
    @{ Cell = 'D7'; Value = '0.8463' }
    @{ Cell = 'E7'; Value = '3.25%' }
    @{ Cell = 'D8'; Value = '0.8601' }
    @{ Cell = 'E8'; Value = '4.00%' }
    @{ Cell = 'D9'; Value = '0.1352' }
    @{ Cell = 'E9'; Value = '1.39%' }
    @{ Cell = 'D10'; Value = '0.06915' }
    @{ Cell = 'E10'; Value = '-0.22%' }
    @{ Cell = 'D11'; Value = '0.02887' }
    @{ Cell = 'E11'; Value = '0.10%' }
    @{ Cell = 'D12'; Value = '0.09377' }
    @{ Cell = 'E12'; Value = '-0.12%' }
    @{ Cell = 'D13'; Value = '0.001513' }
    @{ Cell = 'D14'; Value = '0.04175' }
    @{ Cell = 'E14'; Value = '-9.83%' }
    @{ Cell = 'D15'; Value = '0.0006005' }
    @{ Cell = 'E15'; Value = '-94.03%' }
    @{ Cell = 'D16'; Value = '0.006164' }
    @{ Cell = 'E16'; Value = '-0.68%' }
    @{ Cell = 'E17'; Value = '-4.02%' }
    @{ Cell = 'D18'; Value = '3.033' }
    @{ Cell = 'E18'; Value = '-0.03%' }
    @{ Cell = 'D19'; Value = '2.132' }
    @{ Cell = 'E19'; Value = '-2.31%' }
    @{ Cell = 'E20'; Value = '1.18%' }
    @{ Cell = 'D21'; Value = '0.03348' }
    @{ Cell = 'E21'; Value = '8.36%' }
    @{ Cell = 'E22'; Value = '0.28%' }
    @{ Cell = 'D23'; Value = '3.619' }
    @{ Cell = 'E23'; Value = '-3.43%' }
    @{ Cell = 'E24'; Value = '2.34%' }
    @{ Cell = 'D25'; Value = '0.001211' }
    @{ Cell = 'E25'; Value = '-2.30%' }
    @{ Cell = 'D26'; Value = '0.004444' }
    @{ Cell = 'E26'; Value = '-1.11%' }
    @{ Cell = 'E27'; Value = '22.93%' }
    @{ Cell = 'D28'; Value = '0.0001390' }
    @{ Cell = 'E28'; Value = '-28.23%' }
    @{ Cell = 'D40'; Value = '0.03766' }
    @{ Cell = 'E40'; Value = '3.26%' }
    @{ Cell = 'D41'; Value = '0.005327' }
    @{ Cell = 'E41'; Value = '-13.74%' }
    @{ Cell = 'D42'; Value = '0.1060' }
    @{ Cell = 'E42'; Value = '0.86%' }
    @{ Cell = 'D43'; Value = '0.002288' }
    @{ Cell = 'E43'; Value = '-4.61%' }
    @{ Cell = 'D44'; Value = '0.009283' }
    @{ Cell = 'E44'; Value = '2.97%' }
    @{ Cell = 'D45'; Value = '0.00005102' }
    @{ Cell = 'E45'; Value = '-4.63%' }
    @{ Cell = 'E46'; Value = '-0.01%' }
    @{ Cell = 'D47'; Value = '0.09992' }
    @{ Cell = 'E47'; Value = '-30.58%' }
    @{ Cell = 'D48'; Value = '0.002778' }
    @{ Cell = 'E48'; Value = '19.15%' }
    @{ Cell = 'E49'; Value = '-0.01%' }
    @{ Cell = 'E50'; Value = '-0.01%' }
)

foreach ($u in $numericTextUpdates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
